# Append new "Worldpay recurrent" EUR transaction rows (rows 2-15) to the
# "Transactions" sheet, mirroring the pattern already present in row 1.
# Columns C, D, F, G, H hold numeric-looking reference/amount values that
# must be written as plain text (shared strings), exactly like the
# existing row 1 cells (e.g. "0", "40") already are -- not as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row=2; Cells=@{ "A"="iuliia.4"; "B"="SAR"; "C"="0"; "D"="40"; "E"="MasterCard,null" } },
    @{ Row=3; Cells=@{ "A"="iuliia.4"; "B"="SAR"; "C"="0"; "D"="40"; "E"="MasterCard,null" } },
    @{ Row=4; Cells=@{ "A"="iuliia.4"; "B"="SAR"; "C"="0"; "D"="40"; "E"="MasterCard,null" } },
    @{ Row=5; Cells=@{ "A"="iuliia.4"; "B"="SAR"; "C"="0"; "D"="40"; "E"="MasterCard,Worldpay" } },
    @{ Row=6; Cells=@{ "A"="iuliia.4"; "B"="SAR"; "C"="0"; "D"="40"; "E"="MasterCard,null" } },
    @{ Row=7; Cells=@{ "A"="iuliia.4"; "B"="SAR"; "C"="0"; "D"="40"; "E"="MasterCard,Worldpay"; "F"="3302377965" } },
    @{ Row=8; Cells=@{ "A"="iuliia.4"; "B"="SAR"; "C"="0"; "D"="40"; "E"="MasterCard,null" } },
    @{ Row=9; Cells=@{ "A"="iuliia.4"; "B"="SAR"; "C"="0"; "D"="40"; "E"="MasterCard,null" } },
    @{ Row=10; Cells=@{ "A"="iuliia.4"; "B"="SAR"; "C"="0"; "D"="40"; "E"="MasterCard,Worldpay"; "F"="3340225577" } },
    @{ Row=11; Cells=@{ "A"="iuliia.4"; "B"="SAR"; "C"="0"; "D"="40"; "E"="MasterCard,null" } },
    @{ Row=12; Cells=@{ "A"="iuliia.4"; "B"="SAR"; "C"="0"; "D"="40"; "E"="MasterCard,Worldpay"; "F"="3345760269" } },
    @{ Row=13; Cells=@{ "A"="iuliia.7"; "B"="SAR"; "C"="39.97"; "D"="40"; "E"="MasterCard,null" } },
    @{ Row=14; Cells=@{ "A"="iuliia.7"; "B"="SAR"; "C"="39.97"; "D"="40"; "E"="MasterCard,Worldpay"; "F"="3374872535"; "G"="79.97"; "H"="79.95" } },
    @{ Row=15; Cells=@{ "A"="iuliia.7"; "B"="SAR"; "C"="79.95"; "D"="40"; "E"="MasterCard,Worldpay"; "F"="3389393840"; "G"="119.95"; "H"="119.92" } }
)

foreach ($rowDef in $newRows) {
    $r = $rowDef.Row
    foreach ($col in @("A","B","C","D","E","F","G","H")) {
        if ($rowDef.Cells.ContainsKey($col)) {
            $cellRef = "$col$r"
            $text = $rowDef.Cells[$col]
            $range = $ws.Range($cellRef)

            if ($text -match '^-?[0-9]+(\.[0-9]+)?$') {
                # Force text storage so values like "0", "40", "2016.26"
                # or "3335160798" are kept as shared-string text
                # (matching the existing workbook convention) instead of
                # being auto-detected as numbers. Reset the style back to
                # Normal afterwards so the cell keeps the default
                # (unstyled) appearance.
                $range.NumberFormat = "@"
                $range.Value = $text
                $range.Style = "Normal"
            } else {
                $range.Value = $text
            }
        }
    }
}
